$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftmost "row header" column: shift A1:D1, A2:D2, A3:D3 one
# column to the left, wrapping what used to be in column A into column D.

# Row 1 (numbers 0,1,2 plus the blank header cell)
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = " "

# Row 2 (column headers)
$ws.Range("A2").Value = "Issues traded"
$ws.Range("B2").Value = "Advances"
$ws.Range("C2").Value = "Declines"
$ws.Range("D2").Value = "Scraped @"

# Row 3 (data values) -- keep these as text, not auto-converted numbers.
# Temporarily force a text number format so the comma-containing strings
# aren't parsed as numbers, then restore the default "Normal" style so
# the cell's style index ends up unchanged (as in the original file).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "4,813"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1,696"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "2,811"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "Sep 06, 2022"
